$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C2').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D2').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E2').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F2').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G2').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H2').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B3').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C3').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D3').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E3').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F3').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G3').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H3').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B4').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C4').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D4').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E4').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F4').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G4').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H4').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B5').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C5').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 4.5}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 6.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D5').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E5').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F5').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G5').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H5').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B6').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C6').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D6').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E6').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F6').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G6').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H6').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B7').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C7').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D7').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E7').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F7').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G7').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H7').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B8').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C8').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D8').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 16.5}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 16.5}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E8').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F8').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G8').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H8').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B9').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C9').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 20.4}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D9').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 8.1}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E9').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F9').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G9').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H9').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B10').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C10').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D10').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E10').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F10').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G10').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H10').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B11').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C11').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D11').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E11').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F11').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G11').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H11').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B12').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C12').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D12').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E12').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F12').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G12').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H12').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B13').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C13').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D13').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E13').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F13').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G13').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H13').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B14').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C14').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D14').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E14').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F14').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G14').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H14').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B15').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C15').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D15').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E15').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F15').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G15').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H15').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B16').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C16').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D16').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E16').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F16').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G16').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H16').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B17').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C17').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D17').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E17').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F17').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G17').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H17').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B18').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C18').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D18').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E18').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F18').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G18').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H18').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B19').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 20.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 12.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 20.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C19').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 12.0}, ''D5'': {''第二季'': 20.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 20.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D19').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 20.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 12.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 20.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E19').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 20.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 12.0}, ''D5'': {''第二季'': 20.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F19').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 20.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 12.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 20.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G19').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 20.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 12.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 20.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H19').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 12.0}, ''D5'': {''第二季'': 20.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 20.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B20').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C20').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D20').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E20').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F20').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G20').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H20').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B21').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C21').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D21').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E21').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F21').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G21').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H21').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B22').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 30.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 28.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 20.0}, ''D6'': {''第二季'': 24.0}, ''D7'': {''第二季'': 22.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C22').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 30.0}, ''D2'': {''第二季'': 20.0}, ''D3'': {''第二季'': 28.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 24.0}, ''D7'': {''第二季'': 22.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D22').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 30.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 28.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 20.0}, ''D6'': {''第二季'': 24.0}, ''D7'': {''第二季'': 22.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E22').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 30.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 28.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 24.0}, ''D7'': {''第二季'': 22.0}, ''D8'': {''第二季'': 20.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F22').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 30.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 28.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 20.0}, ''D6'': {''第二季'': 24.0}, ''D7'': {''第二季'': 22.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G22').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 30.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 28.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 20.0}, ''D6'': {''第二季'': 24.0}, ''D7'': {''第二季'': 22.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H22').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 30.0}, ''D2'': {''第二季'': 20.0}, ''D3'': {''第二季'': 28.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 24.0}, ''D7'': {''第二季'': 22.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B23').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C23').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D23').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 18.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 4.5}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E23').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F23').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 21.6}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G23').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 25.8}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H23').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B24').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C24').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D24').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 13.8}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E24').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F24').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 16.5}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G24').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 8.4}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 10.5}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H24').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B25').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C25').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D25').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E25').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F25').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G25').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H25').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B26').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C26').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D26').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E26').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F26').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G26').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H26').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B27').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 1.2}, ''E8'': {''第二季'': 1.2}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 1.2}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 1.2}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 1.2}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C27').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 1.2}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 1.2}, ''E9'': {''第二季'': 1.2}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 1.2}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 1.2}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D27').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 1.2}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 1.2}, ''E8'': {''第二季'': 1.2}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 1.2}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 1.2}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E27').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 1.2}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 1.2}, ''E9'': {''第二季'': 1.2}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 1.2}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 1.2}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F27').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 1.2}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 1.2}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 1.2}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 1.2}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 1.2}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G27').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 1.2}, ''E4'': {''第二季'': 1.2}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 1.2}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 1.2}, ''E13'': {''第二季'': 1.2}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H27').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 1.2}, ''E4'': {''第二季'': 1.2}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 1.2}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 1.2}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 1.2}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B28').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C28').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D28').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E28').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F28').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G28').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H28').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B29').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C29').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D29').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E29').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F29').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G29').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H29').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B30').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C30').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D30').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E30').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F30').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G30').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H30').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B31').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C31').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D31').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E31').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F31').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G31').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H31').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B32').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C32').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D32').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E32').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F32').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G32').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H32').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B33').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C33').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D33').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 13.2}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E33').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F33').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 7.5}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 4.5}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G33').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H33').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B34').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 1.2}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 1.2}, ''E5'': {''第二季'': 1.2}, ''E6'': {''第二季'': 1.2}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 1.2}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C34').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 1.2}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 1.2}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 1.2}, ''E13'': {''第二季'': 1.2}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 1.2}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D34').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 1.2}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 1.2}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 1.2}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 1.2}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 1.2}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E34').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 1.2}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 1.2}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 1.2}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 1.2}, ''E16'': {''第二季'': 1.2}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F34').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 1.2}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 1.2}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 1.2}, ''E10'': {''第二季'': 1.2}, ''E11'': {''第二季'': 1.2}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G34').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 1.2}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 1.2}, ''E11'': {''第二季'': 1.2}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 1.2}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 1.2}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H34').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 1.2}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 1.2}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 1.2}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 1.2}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 1.2}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B35').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 1.2}, ''E3'': {''第二季'': 1.2}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 1.2}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 1.2}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 1.2}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C35').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 1.2}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 1.2}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 1.2}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 1.2}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 1.2}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D35').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 1.2}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 1.2}, ''E6'': {''第二季'': 1.2}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 1.2}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 1.2}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E35').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 1.2}, ''E3'': {''第二季'': 1.2}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 1.2}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 1.2}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 1.2}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F35').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 1.2}, ''E3'': {''第二季'': 1.2}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 1.2}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 1.2}, ''E13'': {''第二季'': 1.2}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G35').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 1.2}, ''E2'': {''第二季'': 1.2}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 1.2}, ''E9'': {''第二季'': 1.2}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 1.2}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H35').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 1.2}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 1.2}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 1.2}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 1.2}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 1.2}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B36').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 1.2}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 1.2}, ''F2'': {''第二季'': 1.2}, ''F3'': {''第二季'': 1.2}, ''F4'': {''第二季'': 1.2}}'
$ws.Range('C36').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 1.2}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 1.2}, ''F2'': {''第二季'': 1.2}, ''F3'': {''第二季'': 1.2}, ''F4'': {''第二季'': 1.2}}'
$ws.Range('D36').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 1.2}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 1.2}, ''F2'': {''第二季'': 1.2}, ''F3'': {''第二季'': 1.2}, ''F4'': {''第二季'': 1.2}}'
$ws.Range('E36').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 1.2}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 1.2}, ''F2'': {''第二季'': 1.2}, ''F3'': {''第二季'': 1.2}, ''F4'': {''第二季'': 1.2}}'
$ws.Range('F36').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 1.2}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 1.2}, ''F2'': {''第二季'': 1.2}, ''F3'': {''第二季'': 1.2}, ''F4'': {''第二季'': 1.2}}'
$ws.Range('G36').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 1.2}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 1.2}, ''F2'': {''第二季'': 1.2}, ''F3'': {''第二季'': 1.2}, ''F4'': {''第二季'': 1.2}}'
$ws.Range('H36').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 1.2}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 1.2}, ''F2'': {''第二季'': 1.2}, ''F3'': {''第二季'': 1.2}, ''F4'': {''第二季'': 1.2}}'
$ws.Range('B37').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C37').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D37').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E37').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F37').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G37').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H37').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B38').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C38').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D38').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E38').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F38').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G38').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H38').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B39').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C39').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D39').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E39').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F39').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G39').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H39').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B40').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C40').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D40').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E40').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F40').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G40').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H40').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B41').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C41').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D41').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E41').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F41').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G41').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H41').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('B42').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('C42').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('D42').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('E42').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('F42').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('G42').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
$ws.Range('H42').Value = '{''A1'': {''第二季'': 0.0}, ''A2'': {''第二季'': 0.0}, ''A3'': {''第二季'': 0.0}, ''A4'': {''第二季'': 0.0}, ''A5'': {''第二季'': 0.0}, ''A6'': {''第二季'': 0.0}, ''B1'': {''第二季'': 0.0}, ''B2'': {''第二季'': 0.0}, ''B3'': {''第二季'': 0.0}, ''B4'': {''第二季'': 0.0}, ''B5'': {''第二季'': 0.0}, ''B6'': {''第二季'': 0.0}, ''B7'': {''第二季'': 0.0}, ''B8'': {''第二季'': 0.0}, ''B9'': {''第二季'': 0.0}, ''B10'': {''第二季'': 0.0}, ''B11'': {''第二季'': 0.0}, ''B12'': {''第二季'': 0.0}, ''B13'': {''第二季'': 0.0}, ''B14'': {''第二季'': 0.0}, ''C1'': {''第二季'': 0.0}, ''C2'': {''第二季'': 0.0}, ''C3'': {''第二季'': 0.0}, ''C4'': {''第二季'': 0.0}, ''C5'': {''第二季'': 0.0}, ''C6'': {''第二季'': 0.0}, ''D1'': {''第二季'': 0.0}, ''D2'': {''第二季'': 0.0}, ''D3'': {''第二季'': 0.0}, ''D4'': {''第二季'': 0.0}, ''D5'': {''第二季'': 0.0}, ''D6'': {''第二季'': 0.0}, ''D7'': {''第二季'': 0.0}, ''D8'': {''第二季'': 0.0}, ''E1'': {''第二季'': 0.0}, ''E2'': {''第二季'': 0.0}, ''E3'': {''第二季'': 0.0}, ''E4'': {''第二季'': 0.0}, ''E5'': {''第二季'': 0.0}, ''E6'': {''第二季'': 0.0}, ''E7'': {''第二季'': 0.0}, ''E8'': {''第二季'': 0.0}, ''E9'': {''第二季'': 0.0}, ''E10'': {''第二季'': 0.0}, ''E11'': {''第二季'': 0.0}, ''E12'': {''第二季'': 0.0}, ''E13'': {''第二季'': 0.0}, ''E14'': {''第二季'': 0.0}, ''E15'': {''第二季'': 0.0}, ''E16'': {''第二季'': 0.0}, ''F1'': {''第二季'': 0.0}, ''F2'': {''第二季'': 0.0}, ''F3'': {''第二季'': 0.0}, ''F4'': {''第二季'': 0.0}}'
